$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (rows 2-10), keeping the header row (row 1) intact.
$ws.Range("A2:T10").Delete()

# Seed the shared-string table so new strings are (re)created in the exact order
# required by the target file: FAPs, MuSCs, Il12a, Il12rb1, ECs.
$ws.Range("A2").Value = "FAPs"
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B2").Value = "Il12a"
$ws.Range("C2").Value = "Il12rb1"
$ws.Range("D2").Value = "ECs"

# Row 2: FAPs | Il12a | Il12rb1 | ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il12a"
$ws.Range("C2").Value = "Il12rb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.338794666666666
$ws.Range("H2").Value = 7.016384
$ws.Range("I2").Value = 0.8670356886266615
$ws.Range("J2").Value = 0.8670356886266616
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06761166666666667
$ws.Range("N2").Value = 0.202835
$ws.Range("O2").Value = 0.04506879582203475
$ws.Range("P2").Value = 0.04506879582203475
$ws.Range("Q2").Value = 0.1581298054044444
$ws.Range("R2").Value = 1.42316824864
$ws.Range("S2").Value = 0.03907625442113231
$ws.Range("T2").Value = 0.03907625442113231

# Row 3: FAPs | Il12a | Il12rb1 | FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il12a"
$ws.Range("C3").Value = "Il12rb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.338794666666666
$ws.Range("H3").Value = 7.016384
$ws.Range("I3").Value = 0.8670356886266615
$ws.Range("J3").Value = 0.8670356886266616
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.211564666666667
$ws.Range("N3").Value = 3.634694
$ws.Range("O3").Value = 0.8076085575052372
$ws.Range("P3").Value = 0.8076085575052372
$ws.Range("Q3").Value = 2.833600980721777
$ws.Range("R3").Value = 25.502408826496
$ws.Range("S3").Value = 0.700225441797338
$ws.Range("T3").Value = 0.7002254417973381

# Row 4: FAPs | Il12a | Il12rb1 | MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Il12a"
$ws.Range("C4").Value = "Il12rb1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.338794666666666
$ws.Range("H4").Value = 7.016384
$ws.Range("I4").Value = 0.8670356886266615
$ws.Range("J4").Value = 0.8670356886266616
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2210116666666666
$ws.Range("N4").Value = 0.6630349999999999
$ws.Range("O4").Value = 0.1473226466727281
$ws.Range("P4").Value = 0.1473226466727281
$ws.Range("Q4").Value = 0.5169009072711109
$ws.Range("R4").Value = 4.65210816544
$ws.Range("S4").Value = 0.1277339924081912
$ws.Range("T4").Value = 0.1277339924081912

# Row 5: MuSCs | Il12a | Il12rb1 | ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Il12a"
$ws.Range("C5").Value = "Il12rb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.358666
$ws.Range("H5").Value = 1.075998
$ws.Range("I5").Value = 0.1329643113733386
$ws.Range("J5").Value = 0.1329643113733386
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06761166666666667
$ws.Range("N5").Value = 0.202835
$ws.Range("O5").Value = 0.04506879582203475
$ws.Range("P5").Value = 0.04506879582203475
$ws.Range("Q5").Value = 0.02425000603666667
$ws.Range("R5").Value = 0.21825005433
$ws.Range("S5").Value = 0.005992541400902449
$ws.Range("T5").Value = 0.005992541400902448

# Row 6: MuSCs | Il12a | Il12rb1 | FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Il12a"
$ws.Range("C6").Value = "Il12rb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.358666
$ws.Range("H6").Value = 1.075998
$ws.Range("I6").Value = 0.1329643113733386
$ws.Range("J6").Value = 0.1329643113733386
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.211564666666667
$ws.Range("N6").Value = 3.634694
$ws.Range("O6").Value = 0.8076085575052372
$ws.Range("P6").Value = 0.8076085575052372
$ws.Range("Q6").Value = 0.4345470527346667
$ws.Range("R6").Value = 3.910923474612
$ws.Range("S6").Value = 0.1073831157078992
$ws.Range("T6").Value = 0.1073831157078992

# Row 7: MuSCs | Il12a | Il12rb1 | MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Il12a"
$ws.Range("C7").Value = "Il12rb1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.358666
$ws.Range("H7").Value = 1.075998
$ws.Range("I7").Value = 0.1329643113733386
$ws.Range("J7").Value = 0.1329643113733386
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2210116666666666
$ws.Range("N7").Value = 0.6630349999999999
$ws.Range("O7").Value = 0.1473226466727281
$ws.Range("P7").Value = 0.1473226466727281
$ws.Range("Q7").Value = 0.07926937043666665
$ws.Range("R7").Value = 0.7134243339299999
$ws.Range("S7").Value = 0.01958865426453696
$ws.Range("T7").Value = 0.01958865426453696
